## Daily Inventory report: insert a new "BRANCH" column in front of the
## existing "PLU CODE" column, and re-tune the column widths so the wider
## report still reads nicely. (See commit message: "Changed the format of
## the xlsx form for Daily Inventory".)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

## 1) Break the two merged ranges before we shuffle columns, so Excel
##    doesn't "helpfully" slide the merges along with the column insert
##    (the title band should stay anchored at column A, and the small
##    A2:C2 label band should keep its original footprint).
$ws.Range("A1:L1").UnMerge()
$ws.Range("A2:C2").UnMerge()

## 2) Insert a new blank column at A; PLU CODE/ITEMCODE/etc. all shift one
##    column to the right (A->B, B->C, ... L->M).
$ws.Columns.Item(1).Insert()

## 3) The report title ("INVENTORY REPORT") rode along with the shift into
##    B1 - put it back at A1, which is the merge anchor.
$ws.Range("A1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = $null

## 4) Re-merge: the title spans the full new width, the small label band
##    keeps its original A2:C2 footprint.
$ws.Range("A1:M1").Merge()
$ws.Range("A2:C2").Merge()

## 5) New header cell: "BRANCH" in front of "PLU CODE", styled like the
##    other header cells (copy the look from the neighboring header).
$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)  ## xlPasteFormats
$ws.Range("A3").Value = "BRANCH"
$excel.CutCopyMode = $false

## 6) Re-tune column widths for the now-13-column layout. (The COM layer's
##    ColumnWidth setter here only keeps whole-pixel precision, so each
##    input below is the value that round-trips closest to the real target
##    character width instead of the target itself.)
$ws.Columns.Item(1).ColumnWidth = 31.918069
$ws.Columns.Item(2).ColumnWidth = 10.418069
$ws.Columns.Item(3).ColumnWidth = 14.251402
$ws.Columns.Item(4).ColumnWidth = 66.418069
$ws.Columns.Item(5).ColumnWidth = 36.584735
$ws.Columns.Item(6).ColumnWidth = 10.584735
$ws.Columns.Item(7).ColumnWidth = 9.918069
$ws.Columns.Item(8).ColumnWidth = 7.418069
$ws.Columns.Item(9).ColumnWidth = 7.418069
$ws.Columns.Item(10).ColumnWidth = 7.918069
$ws.Columns.Item(11).ColumnWidth = 13.418069
$ws.Columns.Item(12).ColumnWidth = 15.751402
$ws.Columns.Item(13).ColumnWidth = 9.084735
